$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "kode_aset" (column B) values are each shifted down by one, i.e. every row now
# shows the code that used to belong to the row above it (row 2 gets a brand
# new "02.01.0001"). Re-write column B, rows 2-12, accordingly.
#
# Row 10's new value ("02.01.0009", no trailing tab) happens to look like a
# date (DD.MM.YYYY) to Excel's auto-detection, so assigning it directly via
# Range.Value would silently convert the cell to a date number. To keep it as
# plain text (matching the original file's habit of storing these codes as
# literal strings with no special number format), build it as a text formula
# in a scratch cell first - string concatenation always yields Text - then
# copy/paste-special just the value into B10 and clean the scratch cell up.

$ws.Range("B2").Value  = "02.01.0001`t"
$ws.Range("B3").Value  = "02.01.0002`t"
$ws.Range("B4").Value  = "02.01.0003`t"
$ws.Range("B5").Value  = "02.01.0004`t"
$ws.Range("B6").Value  = "02.01.0005`t"
$ws.Range("B7").Value  = "02.01.0006`t"
$ws.Range("B8").Value  = "02.01.0007`t"
$ws.Range("B9").Value  = "02.01.0008`t"

$ws.Range("ZZ1").Formula = "=""02.01.000""&""9"""
$ws.Range("ZZ1").Copy()
$ws.Range("B10").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ1").ClearContents()

$ws.Range("B11").Value = "02.01.0010`t"
$ws.Range("B12").Value = "02.01.0011`t"

# The saved workbook's cursor ends up parked one row below the data (B13).
$ws.Range("B13").Select() | Out-Null
